$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "coimbatore"
$ws.Range("A3").Value = "ruralcoimbatore"
$ws.Range("B2").Value = "ambition"
$ws.Range("B3").Value = "ambition"

$ws.Range("C4").Select()
